$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column layout changes -------------------------------------------------
# Column B (start_col) stays the same width but becomes hidden.
$ws.Columns.Item(2).Hidden = $true

# Columns C:E (start_row, end_col, end_row) become zero-width + hidden.
# ColumnWidth always gets +5/6 added internally by this engine, so use a
# small negative input to land exactly on a stored width of 0.
$ws.Range("C1:E1").EntireColumn.ColumnWidth = -0.83
$ws.Range("C1:E1").EntireColumn.Hidden = $true

# --- New column H: crossword answers ("for checking only") -----------------
# Values are written in the same order the shared-string table records them
# (first-use order), so the resulting sharedStrings.xml index assignment
# matches: rows 2-19 top-down, then rows 21-86 top-down, then the header
# (row 1), and finally row 20 (which was corrected last).
$ws.Range("H2").Value = "cass"
$ws.Range("H3").Value = "gulps"
$ws.Range("H4").Value = "shout"
$ws.Range("H5").Value = "lmao"
$ws.Range("H6").Value = "ideal"
$ws.Range("H7").Value = "poets"
$ws.Range("H8").Value = "iota"
$ws.Range("H9").Value = "notmy"
$ws.Range("H10").Value = "rodeo"
$ws.Range("H11").Value = "crime"
$ws.Range("H12").Value = "nap"
$ws.Range("H13").Value = "bic"
$ws.Range("H14").Value = "having"
$ws.Range("H15").Value = "thoughts"
$ws.Range("H16").Value = "ela"
$ws.Range("H17").Value = "div"
$ws.Range("H18").Value = "llcs"
$ws.Range("H19").Value = "wan"
$ws.Range("H21").Value = "bone"
$ws.Range("H22").Value = "coverallthebases"
$ws.Range("H23").Value = "unit"
$ws.Range("H24").Value = "visages"
$ws.Range("H25").Value = "bed"
$ws.Range("H26").Value = "jpeg"
$ws.Range("H27").Value = "tar"
$ws.Range("H28").Value = "pbr"
$ws.Range("H29").Value = "givesthe"
$ws.Range("H30").Value = "degree"
$ws.Range("H31").Value = "ira"
$ws.Range("H32").Value = "tlc"
$ws.Range("H33").Value = "totem"
$ws.Range("H34").Value = "write"
$ws.Range("H35").Value = "about"
$ws.Range("H36").Value = "teri"
$ws.Range("H37").Value = "abram"
$ws.Range("H38").Value = "hipto"
$ws.Range("H39").Value = "yams"
$ws.Range("H40").Value = "silly"
$ws.Range("H41").Value = "steep"
$ws.Range("H42").Value = "ames"
$ws.Range("H43").Value = "cliché"
$ws.Range("H44").Value = "amoral"
$ws.Range("H45").Value = "sativa"
$ws.Range("H46").Value = "soami"
$ws.Range("H47").Value = "gin"
$ws.Range("H48").Value = "udon"
$ws.Range("H49").Value = "letat"
$ws.Range("H50").Value = "pamphlets"
$ws.Range("H51").Value = "sly"
$ws.Range("H52").Value = "sprigs"
$ws.Range("H53").Value = "hooch"
$ws.Range("H54").Value = "oed"
$ws.Range("H55").Value = "ute"
$ws.Range("H56").Value = "tso"
$ws.Range("H57").Value = "endor"
$ws.Range("H58").Value = "buc"
$ws.Range("H59").Value = "gina"
$ws.Range("H60").Value = "oldhat"
$ws.Range("H61").Value = "twos"
$ws.Range("H62").Value = "sane"
$ws.Range("H63").Value = "velvet"
$ws.Range("H64").Value = "nes"
$ws.Range("H65").Value = "jet"
$ws.Range("H66").Value = "slightbit"
$ws.Range("H67").Value = "bas"
$ws.Range("H68").Value = "cub"
$ws.Range("H69").Value = "oneg"
$ws.Range("H70").Value = "vidi"
$ws.Range("H71").Value = "egad"
$ws.Range("H72").Value = "beret"
$ws.Range("H73").Value = "jeremy"
$ws.Range("H74").Value = "psa"
$ws.Range("H75").Value = "prteam"
$ws.Range("H76").Value = "beerme"
$ws.Range("H77").Value = "remiss"
$ws.Range("H78").Value = "vital"
$ws.Range("H79").Value = "elope"
$ws.Range("H80").Value = "gotya"
$ws.Range("H81").Value = "cute"
$ws.Range("H82").Value = "was"
$ws.Range("H83").Value = "rbi"
$ws.Range("H84").Value = "irl"
$ws.Range("H85").Value = "ahs"
$ws.Range("H86").Value = "top"
$ws.Range("H1").Value = "answer (optional column, for checking only)"
$ws.Range("H20").Value = "jonesed"

# --- Selection: last thing the author clicked on before saving -------------
$ws.Range("H21").Select()
